# "cleaning up last flow"
# Updates a handful of stat figures on Sheet1/Sheet2, tidies the duplicate
# underline-font cell style that had crept onto B12/B13/B15 of Sheet1 (now
# matching the rest of column B), and leaves the selection/active-sheet
# where the author last left it (Sheet2 active, with Sheet1 scrolled down a
# touch and Sheet2's cursor sitting on D6).

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# ---------------------------------------------------------------------
# Sheet1 value corrections
# ---------------------------------------------------------------------
$ws1.Range("B2").Value = 722
$ws1.Range("C2").Value = 168
$ws1.Range("D2").Value = 684
$ws1.Range("E2").Value = 95
$ws1.Range("F2").Value = 0.247
$ws1.Range("G2").Value = 0.755

$ws1.Range("B3").Value = 644

$ws1.Range("B5").Value = 741
$ws1.Range("C5").Value = 207
$ws1.Range("E5").Value = 148
$ws1.Range("F5").Value = 0.274
$ws1.Range("G5").Value = 0.788

$ws1.Range("F10").Value = 0.26
$ws1.Range("G10").Value = 0.77

$ws1.Range("B12").Value = 688
$ws1.Range("B13").Value = 688
$ws1.Range("C13").Value = 169

# B12/B13/B15 were carrying a leftover underline-font style (a near-dupe of
# the font already used elsewhere in column B). Re-apply the normal column
# style B14 uses so they fall back in line with the rest of the column.
$ws1.Range("B14").Copy()
$ws1.Range("B12").PasteSpecial(-4122)
$ws1.Range("B13").PasteSpecial(-4122)
$ws1.Range("B15").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------
# Sheet2 value corrections
# ---------------------------------------------------------------------
$ws2.Range("D2").Value = 3.9
$ws2.Range("E2").Value = 1.2

$ws2.Range("D6").Value = 4.65
$ws2.Range("E6").Value = 1.18

$ws2.Range("D8").Value = 3.48
$ws2.Range("E8").Value = 1.11

# ---------------------------------------------------------------------
# Cursor / active sheet bookkeeping
# ---------------------------------------------------------------------
$ws1.Range("F16").Select() | Out-Null
$ws2.Range("D6").Select() | Out-Null
$ws2.Activate()
